# Update lecture 7. Add Assignment 9
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Assignment 9" column header ---
$ws.Range("L1").Value = "Assignment 9"

# --- Bold the whole header row (new font + cell style for A1:L1) ---
$ws.Range("A1:L1").Font.Bold = $true

# --- Lecture 7 grade updates (new participation / assignment marks) ---
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1

$ws.Range("I4").Value = 1

$ws.Range("C6").Value = 4

$ws.Range("K7").Value = 1

$ws.Range("K9").Value = 1

$ws.Range("J12").Value = 1

$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 1

$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 1

$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1

$ws.Range("J17").Value = 1

$ws.Range("J18").Value = 1

$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 1

$ws.Range("J21").Value = 1

$ws.Range("K23").Value = 1

$ws.Range("J25").Value = 1
$ws.Range("K25").Value = 1

$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 1

$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 1

$ws.Range("K30").Value = 1

# --- Move the active selection (matches the saved view state) ---
[void]$ws.Range("N23").Select()
